$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 39: turn the previously-blank row into a real data row, using
# the exact same per-column formatting as row 38 (the last data row).
# ------------------------------------------------------------------
$ws.Range("A38:O38").Copy()
$ws.Range("A39:O39").PasteSpecial(-4122)

$ws.Range("A39").Value = $ws.Range("A38").Value2
$ws.Range("B39").Value = "14006 x 1081"
$ws.Range("C39").Value = $ws.Range("C38").Value2
$ws.Range("D39").Value = $ws.Range("D38").Value2
$ws.Range("E39").Value = $ws.Range("E38").Value2
$ws.Range("F39").Value = $ws.Range("F38").Value2
$ws.Range("G39").Value = $ws.Range("G38").Value2
$ws.Range("H39").Value = $ws.Range("H38").Value2
$ws.Range("I39").Value = 50
$ws.Range("K39").Value = 115
$ws.Range("L39").Value = "84.2 & 72.5"
$ws.Range("M39").Value = "79.7 & 71.0"
$ws.Range("N39").Value = 20
$ws.Range("O39").Value = 62.8

# ------------------------------------------------------------------
# Rows 41-49: nine more trailing blank rows, formatted like the old
# trailing blank row 40 (copy its look before we touch row 40 itself).
# ------------------------------------------------------------------
$ws.Range("A40:O40").Copy()
for ($r = 41; $r -le 49; $r++) {
    $ws.Range(("A{0}:O{0}" -f $r)).PasteSpecial(-4122) | Out-Null
}

# ------------------------------------------------------------------
# Rows 40-48 become short "spacer" rows: slightly shorter and with a
# lighter (automatic/theme) font colour on the numeric columns, while
# row 49 keeps the original look (height 19.5, solid black numbers).
# ------------------------------------------------------------------
for ($r = 40; $r -le 48; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
    $ws.Range("I" + $r).Font.ThemeColor = 1
    $ws.Range("K" + $r).Font.ThemeColor = 1
    $ws.Range("N" + $r).Font.ThemeColor = 1
    $ws.Range("O" + $r).Font.ThemeColor = 1
}
$ws.Rows.Item(49).RowHeight = 19.5

Write-Output "done"
